# Adds 5 new quarterly columns (D:H) in front of the existing quarterly data,
# shifting the previous D:H data right to I:M, and fills the new columns with
# the earlier-period figures (cumulative-data / predict_income refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 new blank columns at D:H - this shifts existing D:H (and their
#    number formats / styles / shared-string references) to I:M automatically.
$ws.Range("D1:H1").EntireColumn.Insert()

# 2) Re-apply the column widths exactly as in the target layout.
#    (ColumnWidth set via COM is offset ~0.83 from the stored OOXML width,
#    so subtract that to land on the exact target width.)
$ws.Range("D1").EntireColumn.ColumnWidth = 28.17   # -> stored width 29
$ws.Range("E1").EntireColumn.ColumnWidth = 28.17   # -> stored width 29
$ws.Range("F1").EntireColumn.ColumnWidth = 30.17   # -> stored width 31
$ws.Range("G1").EntireColumn.ColumnWidth = 28.17   # -> stored width 29
$ws.Range("H1").EntireColumn.ColumnWidth = 28.17   # -> stored width 29

# 3) Header row 8: quarterly period labels for the newly inserted columns.
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

# 4) Row 9: publish dates for the newly inserted columns.
$ws.Range("D9").Value = "1400-10-09 (5)"
$ws.Range("E9").Value = "1400-11-24 (3)"
$ws.Range("F9").Value = "1401-04-21 (11)"
$ws.Range("G9").Value = "1401-04-29 (2)"
$ws.Range("H9").Value = "1401-09-15 (6)"

# 5) Income-statement figures for the newly inserted columns (rows 11-27).
$rowValues = @{
    11 = @(4369126, 7166921, 3287713, 5705755, 5830143)
    12 = @(-1474985, -2269924, -1703169, -2071476, -2513373)
    13 = @(2894141, 4896997, 1584544, 3634279, 3316770)
    14 = @(-169483, -206784, -242555, -325171, -412071)
    15 = @(0, 0, 0, 0, 0)
    16 = @(904670, -83613, 1776268, 31440, 1180597)
    17 = @(3629328, 4606600, 3118257, 3340548, 4085296)
    18 = @(-41858, -26275, -51686, -12094, -46126)
    19 = @(42922, 76131, 76490, 36688, 43585)
    20 = @(3630392, 4656456, 3143061, 3365142, 4082755)
    21 = @(-271472, -180163, -19099, -265375, -137522)
    22 = @(3358920, 4476293, 3123962, 3099767, 3945233)
    23 = @(0, 0, 0, 0, 0)
    24 = @(3358920, 4476293, 3123962, 3099767, 3945233)
    25 = @(560, 746, 521, 517, 658)
    26 = @(6000000, 6000000, 6000000, 6000000, 6000000)
    27 = @(560, 746, 521, 517, 658)
}

$cols = @("D", "E", "F", "G", "H")
foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
